# Drop the redundant "Q" label cells above each long question - the
# question text itself (in column B) already conveys that it's the
# question, so the separate one-cell "Q" row is just noise.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").ClearContents() | Out-Null
$ws.Range("A10").ClearContents() | Out-Null

# Park the selection on the (now empty) A10 cell.
$ws.Range("A10").Select() | Out-Null
